$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6193.8335
$ws.Range("I19").Value = 164.66667
$ws.Range("J19").Value = 9811.333000000001
$ws.Range("K19").Value = 164.66667
$ws.Range("L19").Value = 9811.333000000001
$ws.Range("M19").Value = 10.33332999999999
$ws.Range("N19").Value = -10161.333
$ws.Range("H137").Value = 2943956.2
$ws.Range("I137").Value = 3848835
$ws.Range("J137").Value = 3100
$ws.Range("K137").Value = 11546505
$ws.Range("L137").Value = 9300
$ws.Range("M137").Value = -11543955
$ws.Range("N137").Value = -14400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3765.5244
$ws.Range("I32").Value = 2639.7576
$ws.Range("J32").Value = 8409.3125
$ws.Range("K32").Value = 2639.7576
$ws.Range("L32").Value = 8409.3125
$ws.Range("M32").Value = -2352.7576
$ws.Range("N32").Value = -8983.3125
$ws.Range("H52").Value = 39780
$ws.Range("J52").Value = 39780
$ws.Range("L52").Value = 39780
$ws.Range("N52").Value = -40416
$ws.Range("H61").Value = 3124.1333
$ws.Range("I61").Value = 1474.9333
$ws.Range("J61").Value = 4773.3335
$ws.Range("K61").Value = 1474.9333
$ws.Range("L61").Value = 4773.3335
$ws.Range("M61").Value = -1262.9333
$ws.Range("N61").Value = -5197.3335
$ws.Range("H63").Value = 5163.636
$ws.Range("J63").Value = 8400
$ws.Range("L63").Value = 8400
$ws.Range("N63").Value = -9772
$ws.Range("H66").Value = 5163.636
$ws.Range("J66").Value = 8400
$ws.Range("L66").Value = 42000
$ws.Range("N66").Value = -48864
$ws.Range("H88").Value = 4111.2856
$ws.Range("I88").Value = 1995
$ws.Range("J88").Value = 6933
$ws.Range("K88").Value = 1995
$ws.Range("L88").Value = 6933
$ws.Range("M88").Value = -1589
$ws.Range("N88").Value = -7745
$ws.Range("H91").Value = 4111.2856
$ws.Range("I91").Value = 1995
$ws.Range("J91").Value = 6933
$ws.Range("K91").Value = 1995
$ws.Range("L91").Value = 6933
$ws.Range("M91").Value = -591
$ws.Range("N91").Value = -9741
$ws.Range("H110").Value = 1382.3103
$ws.Range("I110").Value = 708.43475
$ws.Range("J110").Value = 3965.5
$ws.Range("K110").Value = 708.43475
$ws.Range("L110").Value = 3965.5
$ws.Range("M110").Value = 1336.56525
$ws.Range("N110").Value = -8055.5
$ws.Range("H136").Value = 3124.1333
$ws.Range("I136").Value = 1474.9333
$ws.Range("J136").Value = 4773.3335
$ws.Range("K136").Value = 4424.7999
$ws.Range("L136").Value = 14320.0005
$ws.Range("M136").Value = -1874.7999
$ws.Range("N136").Value = -19420.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1388.08
$ws.Range("I107").Value = 969.44446
$ws.Range("J107").Value = 2464.5715
$ws.Range("K107").Value = 969.44446
$ws.Range("L107").Value = 2464.5715
$ws.Range("M107").Value = 950.55554
$ws.Range("N107").Value = -6304.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2634825.2
$ws.Range("I31").Value = 5001914
$ws.Range("J31").Value = 4727.1113
$ws.Range("K31").Value = 5001914
$ws.Range("L31").Value = 4727.1113
$ws.Range("M31").Value = -5001619
$ws.Range("N31").Value = -5317.1113
$ws.Range("H34").Value = 2634825.2
$ws.Range("I34").Value = 5001914
$ws.Range("J34").Value = 4727.1113
$ws.Range("K34").Value = 5001914
$ws.Range("L34").Value = 4727.1113
$ws.Range("M34").Value = -5001712
$ws.Range("N34").Value = -5131.1113
$ws.Range("H58").Value = 22730624
$ws.Range("I58").Value = 2087.2
$ws.Range("J58").Value = 41671070
$ws.Range("K58").Value = 2087.2
$ws.Range("L58").Value = 41671070
$ws.Range("M58").Value = -1884.2
$ws.Range("N58").Value = -41671476
$ws.Range("H86").Value = 3743.8518
$ws.Range("I86").Value = 3109.7273
$ws.Range("J86").Value = 4179.8125
$ws.Range("K86").Value = 3109.7273
$ws.Range("L86").Value = 4179.8125
$ws.Range("M86").Value = -1986.7273
$ws.Range("N86").Value = -6425.8125
$ws.Range("H89").Value = 3743.8518
$ws.Range("I89").Value = 3109.7273
$ws.Range("J89").Value = 4179.8125
$ws.Range("K89").Value = 15548.6365
$ws.Range("L89").Value = 20899.0625
$ws.Range("M89").Value = -9932.636500000001
$ws.Range("N89").Value = -32131.0625
$ws.Range("H136").Value = 22730624
$ws.Range("I136").Value = 2087.2
$ws.Range("J136").Value = 41671070
$ws.Range("K136").Value = 6261.599999999999
$ws.Range("L136").Value = 125013210
$ws.Range("M136").Value = -3711.599999999999
$ws.Range("N136").Value = -125018310

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2433.3333
$ws.Range("J64").Value = 3000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9540
$ws.Range("H67").Value = 2433.3333
$ws.Range("J67").Value = 3000
$ws.Range("L67").Value = 9000
$ws.Range("N67").Value = -10872
$ws.Range("H68").Value = 2365.9556
$ws.Range("I68").Value = 703.6667
$ws.Range("J68").Value = 3820.4583
$ws.Range("K68").Value = 2111.0001
$ws.Range("L68").Value = 11461.3749
$ws.Range("M68").Value = -1300.0001
$ws.Range("N68").Value = -13083.3749
$ws.Range("H71").Value = 2365.9556
$ws.Range("I71").Value = 703.6667
$ws.Range("J71").Value = 3820.4583
$ws.Range("K71").Value = 6333.0003
$ws.Range("L71").Value = 34384.1247
$ws.Range("M71").Value = -2277.0003
$ws.Range("N71").Value = -42496.1247
$ws.Range("H74").Value = 3422.2222
$ws.Range("I74").Value = 2333.3333
$ws.Range("J74").Value = 3966.6667
$ws.Range("K74").Value = 6999.999899999999
$ws.Range("L74").Value = 11900.0001
$ws.Range("M74").Value = -5938.999899999999
$ws.Range("N74").Value = -14022.0001
$ws.Range("H77").Value = 3422.2222
$ws.Range("I77").Value = 2333.3333
$ws.Range("J77").Value = 3966.6667
$ws.Range("K77").Value = 20999.9997
$ws.Range("L77").Value = 35700.0003
$ws.Range("M77").Value = -15695.9997
$ws.Range("N77").Value = -46308.0003
$ws.Range("H92").Value = 3050.5
$ws.Range("J92").Value = 3733.3333
$ws.Range("L92").Value = 11199.9999
$ws.Range("N92").Value = -13695.9999
$ws.Range("H137").Value = 3464.1538
$ws.Range("I137").Value = 3186.4119
$ws.Range("J137").Value = 3988.7778
$ws.Range("K137").Value = 9559.235700000001
$ws.Range("L137").Value = 11966.3334
$ws.Range("M137").Value = -4459.235700000001
$ws.Range("N137").Value = -22166.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 243003.6
$ws.Range("J18").Value = 53754.5
$ws.Range("L18").Value = 53754.5
$ws.Range("N18").Value = -54340.5
$ws.Range("H132").Value = 3218.2703
$ws.Range("I132").Value = 2245.7036
$ws.Range("J132").Value = 5844.2
$ws.Range("K132").Value = 6737.110799999999
$ws.Range("L132").Value = 17532.6
$ws.Range("M132").Value = -4207.110799999999
$ws.Range("N132").Value = -22592.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 125029500
$ws.Range("I4").Value = 200014190
$ws.Range("J4").Value = 55006.668
$ws.Range("K4").Value = 200014190
$ws.Range("L4").Value = 55006.668
$ws.Range("M4").Value = -200014077
$ws.Range("N4").Value = -55232.668
$ws.Range("H5").Value = 143290
$ws.Range("J5").Value = 166671.67
$ws.Range("L5").Value = 166671.67
$ws.Range("N5").Value = -166897.67
$ws.Range("H14").Value = 236819.45
$ws.Range("I14").Value = 419502
$ws.Range("J14").Value = 17600.4
$ws.Range("K14").Value = 419502
$ws.Range("L14").Value = 17600.4
$ws.Range("M14").Value = -419330
$ws.Range("N14").Value = -17944.4
$ws.Range("H28").Value = 125029500
$ws.Range("I28").Value = 200014190
$ws.Range("J28").Value = 55006.668
$ws.Range("K28").Value = 200014190
$ws.Range("L28").Value = 55006.668
$ws.Range("M28").Value = -200013958
$ws.Range("N28").Value = -55470.668
$ws.Range("H37").Value = 125029500
$ws.Range("I37").Value = 200014190
$ws.Range("J37").Value = 55006.668
$ws.Range("K37").Value = 200014190
$ws.Range("L37").Value = 55006.668
$ws.Range("M37").Value = -200014083
$ws.Range("N37").Value = -55220.668
$ws.Range("H40").Value = 2228.7646
$ws.Range("I40").Value = 1617.8
$ws.Range("J40").Value = 2483.3333
$ws.Range("K40").Value = 1617.8
$ws.Range("L40").Value = 2483.3333
$ws.Range("M40").Value = -1481.8
$ws.Range("N40").Value = -2755.3333
$ws.Range("H55").Value = 2976.5881
$ws.Range("I55").Value = 2168.182
$ws.Range("J55").Value = 4458.6665
$ws.Range("K55").Value = 2168.182
$ws.Range("L55").Value = 4458.6665
$ws.Range("M55").Value = -1995.182
$ws.Range("N55").Value = -4804.6665
$ws.Range("H104").Value = 19181.428
$ws.Range("J104").Value = 19181.428
$ws.Range("L104").Value = 19181.428
$ws.Range("N104").Value = -26169.428
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -35060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2683.3333
$ws.Range("J113").Value = 3025
$ws.Range("L113").Value = 9075
$ws.Range("N113").Value = -13415
$ws.Range("H132").Value = 339552.8
$ws.Range("I132").Value = 557828.4
$ws.Range("K132").Value = 1673485.2
$ws.Range("M132").Value = -1670955.2
